$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "69.477.08"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "3.393.16"
$ws.Range("E3").Value = "  +4.83%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCell "D5" "191.62"
$ws.Range("E5").Value = "  +4.19%  "
Set-TextCell "D6" "593.66"
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("E7").Value = "  +0.04%  "
Set-TextCell "D8" "0.607"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").Value = "3.983.60"
$ws.Range("E12").Value = "  +4.91%  "
Set-TextCell "D13" "0.138"
$ws.Range("E13").Value = "  +0.64%  "
Set-TextCell "D14" "28.73"
$ws.Range("E14").Value = "  +3.94%  "
$ws.Range("D15").Value = "69.531.77"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "3.370.66"
$ws.Range("E17").Value = "  +3.74%  "
Set-TextCell "D18" "449.93"
$ws.Range("E18").Value = "  +13.83%  "
Set-TextCell "D19" "5.85"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("E20").Value = "  +2.61%  "
Set-TextCell "D21" "7.83"
$ws.Range("E21").Value = "  +3.62%  "
Set-TextCell "D22" "75.52"
$ws.Range("E22").Value = "  +6.02%  "
$ws.Range("D24").Value = "3.525.91"
$ws.Range("E24").Value = "  +4.43%  "
$ws.Range("E25").Value = "  +4.32%  "
Set-TextCell "D26" "0.523"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  +0.29%  "
Set-TextCell "D30" "2.00"
$ws.Range("E30").Value = "  +2.31%  "
Set-TextCell "D31" "23.51"
$ws.Range("E31").Value = "  +3.96%  "
Set-TextCell "D32" "5.68"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  +3.29%  "
Set-TextCell "D34" "7.00"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +0.02%  "
Set-TextCell "D36" "1.57"
$ws.Range("E36").Value = "  +6.74%  "
Set-TextCell "D37" "164.79"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("E38").Value = "  +3.49%  "
Set-TextCell "D39" "27.96"
$ws.Range("E39").Value = "  +5.58%  "
Set-TextCell "D40" "0.819"
$ws.Range("E40").Value = "  +1.95%  "
Set-TextCell "D41" "4.61"
$ws.Range("E41").Value = "  +1.46%  "
Set-TextCell "D42" "6.63"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("D43").Value = "2.754.35"
$ws.Range("E43").Value = "  +5.37%  "
$ws.Range("E44").Value = "  +2.56%  "
Set-TextCell "D45" "25.57"
$ws.Range("E45").Value = "  +3.68%  "
Set-TextCell "D46" "0.0692"
$ws.Range("E46").Value = "  +0.74%  "
Set-TextCell "D47" "40.96"
$ws.Range("E47").Value = "  +1.02%  "
Set-TextCell "D48" "340.35"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E49").Value = "  +2.38%  "
Set-TextCell "D50" "33.05"
$ws.Range("E50").Value = "  +7.54%  "
$ws.Range("E51").Value = "  +6.12%  "
